$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells B1:G1 from A,B,C,D,E,F to rating1..rating6
$ws.Range("B1").Value = "rating1"
$ws.Range("C1").Value = "rating2"
$ws.Range("D1").Value = "rating3"
$ws.Range("E1").Value = "rating4"
$ws.Range("F1").Value = "rating5"
$ws.Range("G1").Value = "rating6"

# Select B1:G1 as shown in the diff
$ws.Range("B1:G1").Select()

# Apply AutoFilter over the used data range
$ws.Range("A1:G58").AutoFilter()

# Register the hidden sheet-scoped _FilterDatabase defined name that Excel
# creates automatically when AutoFilter is turned on
$fd = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:G58"))
$fd.Visible = $false
